$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "TSE-Overall-Index-TEDPIX"
$ws.Range("B10").Value = "TSE Overall Index (TEDPIX)"
$ws.Range("C10").Value = "https://github.com/imahdimir/d-TSE-Overall-Index-TEDPIX"
$ws.Range("D10").Value = "https://raw.github.com/imahdimir/d-TSE-Overall-Index-TEDPIX/main/META.json"

# Row 11
$ws.Range("A11").Value = "TSE-working-days"
$ws.Range("B11").Value = "TSE Working Days"
$ws.Range("C11").Value = "https://github.com/imahdimir/d-TSE-working-days"
$ws.Range("D11").Value = "https://raw.github.com/imahdimir/d-TSE-working-days/main/META.json"

# Row 12
$ws.Range("A12").Value = "Ticker-2-BaseTicker"
$ws.Range("B12").Value = "Tickers to BaseTickers map"
$ws.Range("C12").Value = "https://github.com/imahdimir/d-Ticker-2-BaseTicker"
$ws.Range("D12").Value = "https://raw.github.com/imahdimir/d-Ticker-2-BaseTicker/main/META.json"

# Row 13
$ws.Range("A13").Value = "FirmTicker-IPO_JDate"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "https://github.com/imahdimir/d-FirmTicker-IPO_JDate"
$ws.Range("D13").Value = "https://raw.github.com/imahdimir/d-FirmTicker-IPO_JDate/main/META.json"

# Row 14
$ws.Range("A14").Value = "firms-adjusted-Prices-1-OHLCL-daily"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "https://github.com/imahdimir/d-firms-adjusted-Prices-1-OHLCL-daily"
$ws.Range("D14").Value = "https://raw.github.com/imahdimir/d-firms-adjusted-Prices-1-OHLCL-daily/main/META.json"

# Row 15
$ws.Range("A15").Value = "FirmTicker-possible-trade-spells"
$ws.Range("C15").Value = "https://github.com/imahdimir/d-FirmTicker-possible-trade-spells"
$ws.Range("D15").Value = "https://raw.github.com/imahdimir/d-FirmTicker-possible-trade-spells/main/META.json"

# Row 16
$ws.Range("A16").Value = "FirmTicker-open-duration-daily"
$ws.Range("C16").Value = "https://github.com/imahdimir/d-FirmTicker-open-duration-daily"
$ws.Range("D16").Value = "https://raw.github.com/imahdimir/d-FirmTicker-open-duration-daily/main/META.json"

# Row 17
$ws.Range("A17").Value = "CodalTicker-2-ISIC"
$ws.Range("C17").Value = "https://github.com/imahdimir/d-CodalTicker-2-ISIC"
$ws.Range("D17").Value = "https://raw.github.com/imahdimir/d-CodalTicker-2-ISIC/main/META.json"

# Row 18
$ws.Range("A18").Value = "USD-IRR-monthly"
$ws.Range("C18").Value = "https://github.com/imahdimir/d-USD-IRR-monthly"
$ws.Range("D18").Value = "https://raw.github.com/imahdimir/d-USD-IRR-monthly/main/META.json"

# Row 19
$ws.Range("A19").Value = "FirmTicker-2-ISIC"
$ws.Range("C19").Value = "https://github.com/imahdimir/d-FirmTicker-2-ISIC"
$ws.Range("D19").Value = "https://raw.github.com/imahdimir/d-FirmTicker-2-ISIC/main/META.json"

# Row 20
$ws.Range("A20").Value = "FirmTicker-Industry-SubIndustry"
$ws.Range("C20").Value = "https://github.com/imahdimir/d-FirmTicker-Industry-SubIndustry"
$ws.Range("D20").Value = "https://raw.github.com/imahdimir/d-FirmTicker-Industry-SubIndustry/main/META.json"

# Row 21 (new)
$ws.Range("A21").Value = "Ticker-2-FirmTicker"
$ws.Range("C21").Value = "https://github.com/imahdimir/d-Ticker-2-FirmTicker"
$ws.Range("D21").Value = "https://raw.github.com/imahdimir/d-Ticker-2-FirmTicker/main/META.json"

# Row 22 (new)
$ws.Range("A22").Value = "FirmTicker-status-change"
$ws.Range("B22").Value = "Cleaned Status changes of each TSETMC_ID on TESTMC.com"
$ws.Range("C22").Value = "https://github.com/imahdimir/d-FirmTicker-status-change"
$ws.Range("D22").Value = "https://raw.github.com/imahdimir/d-FirmTicker-status-change/main/META.json"

# Row 23 (new)
$ws.Range("A23").Value = "TSETMC_ID-Shenase"
$ws.Range("C23").Value = "https://github.com/imahdimir/d-TSETMC_ID-Shenase"
$ws.Range("D23").Value = "https://raw.github.com/imahdimir/d-TSETMC_ID-Shenase/main/META.json"

# Row 24 (new)
$ws.Range("A24").Value = "IFB-stocks-only-TradeValue-monthly"
$ws.Range("C24").Value = "https://github.com/imahdimir/d-IFB-stocks-only-TradeValue-monthly"
$ws.Range("D24").Value = "https://raw.github.com/imahdimir/d-IFB-stocks-only-TradeValue-monthly/main/META.json"

# Row 25 (new)
$ws.Range("A25").Value = "FirmTicker-DPS"
$ws.Range("C25").Value = "https://github.com/imahdimir/d-FirmTicker-DPS"
$ws.Range("D25").Value = "https://raw.github.com/imahdimir/d-FirmTicker-DPS/main/META.json"

# Row 26 (new)
$ws.Range("A26").Value = "CompanyName-2-FirmTicker"
$ws.Range("C26").Value = "https://github.com/imahdimir/d-CompanyName-2-FirmTicker"
$ws.Range("D26").Value = "https://raw.github.com/imahdimir/d-CompanyName-2-FirmTicker/main/META.json"

# Row 27 (new)
$ws.Range("A27").Value = "FirmTickers"
$ws.Range("C27").Value = "https://github.com/imahdimir/d-FirmTickers"
$ws.Range("D27").Value = "https://raw.github.com/imahdimir/d-FirmTickers/main/META.json"
